$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row with machine-readable column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix capitalization of "de"/"del"/"de la" -> "De"/"De La" in a handful of place names
$ws.Range("B2").Value = "Pabellón De Arteaga"
$ws.Range("B8").Value = "Mazapa De Madero"
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("A16").Value = "Estado De México"
$ws.Range("B16").Value = "Ecatepec De Morelos"
$ws.Range("B18").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B21").Value = "Acapulco De Juárez"

# Remove trailing metadata/footer rows 47-51 (sample size, source, author, date)
$ws.Range("A47:A51").EntireRow.Delete()
